$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'295.34"
$ws.Range("E2").Value = "'-4.47%"
$ws.Range("D3").Value = "'40.13"
$ws.Range("E3").Value = "'-2.65%"
$ws.Range("D4").Value = "'5.022"
$ws.Range("E4").Value = "'-3.68%"
$ws.Range("D5").Value = "'0.07380"
$ws.Range("E5").Value = "'-4.01%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.529"
$ws.Range("E6").Value = "'-7.04%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9214"
$ws.Range("E7").Value = "'0.63%"
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").Value = "'0.1170"
$ws.Range("E8").Value = "'-5.95%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1756"
$ws.Range("E9").Value = "'-3.90%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.08620"
$ws.Range("E10").Value = "'-5.91%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.04167"
$ws.Range("E11").Value = "'-1.43%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.1054"
$ws.Range("E12").Value = "'0.23%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001271"
$ws.Range("E13").Value = "'1.35%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "'0.005765"
$ws.Range("E14").Value = "'-0.82%"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").Value = "'3.374"
$ws.Range("E15").Value = "'0.74%"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").Value = "'4.314"
$ws.Range("E16").Value = "'-0.05%"
$ws.Range("D19").Value = "'7.589"
$ws.Range("E19").Value = "'2.72%"
$ws.Range("E20").Value = "'-3.28%"
$ws.Range("D21").Value = "'0.2811"
$ws.Range("E21").Value = "'-0.36%"
$ws.Range("D22").Value = "'0.03825"
$ws.Range("E22").Value = "'-5.24%"
$ws.Range("D23").Value = "'0.001285"
$ws.Range("E23").Value = "'1.45%"
$ws.Range("D24").Value = "'0.003898"
$ws.Range("E24").Value = "'-4.87%"
$ws.Range("E25").Value = "'0.22%"
$ws.Range("D26").Value = "'0.0003735"
$ws.Range("E26").Value = "'-95.03%"
$ws.Range("D38").Value = "'0.02303"
$ws.Range("E38").Value = "'-9.99%"
$ws.Range("D39").Value = "'0.05003"
$ws.Range("E39").Value = "'-6.49%"
$ws.Range("D40").Value = "'0.007709"
$ws.Range("E40").Value = "'-1.55%"
$ws.Range("D41").Value = "'0.1273"
$ws.Range("E41").Value = "'-3.23%"
$ws.Range("E42").Value = "'115.64%"
$ws.Range("D43").Value = "'0.007414"
$ws.Range("E43").Value = "'11.20%"
$ws.Range("D44").Value = "'0.007687"
$ws.Range("E44").Value = "'-4.68%"
$ws.Range("D45").Value = "'0.3207"
$ws.Range("E45").Value = "'4.46%"
$ws.Range("D46").Value = "'0.00006465"
$ws.Range("E46").Value = "'-4.08%"
$ws.Range("E47").Value = "'0.20%"
$ws.Range("E48").Value = "'6.98%"
$ws.Range("D49").Value = "'0.004215"
$ws.Range("E49").Value = "'35.81%"
$ws.Range("D50").Value = "'0.00002107"
$ws.Range("E50").Value = "'0.20%"
$ws.Range("D51").Value = "'0.0002007"
$ws.Range("E51").Value = "'0.20%"
